$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 32.93949833333333
$ws.Range("H2").Value = 98.818495
$ws.Range("I2").Value = 0.02571831923682078
$ws.Range("J2").Value = 0.02571831923682077
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2262196666666666
$ws.Range("N2").Value = 0.6786589999999999
$ws.Range("O2").Value = 0.03145179203784564
$ws.Range("P2").Value = 0.03145179203784564
$ws.Range("Q2").Value = 7.451562333133888
$ws.Range("R2").Value = 67.06406099820499
$ws.Range("S2").Value = 0.0008088872281994122
$ws.Range("T2").Value = 0.0008088872281994121

# Row 3
$ws.Range("G3").Value = 32.93949833333333
$ws.Range("H3").Value = 98.818495
$ws.Range("I3").Value = 0.02571831923682078
$ws.Range("J3").Value = 0.02571831923682077
$ws.Range("O3").Value = 0.9636438974901603
$ws.Range("P3").Value = 0.9636438974901604
$ws.Range("Q3").Value = 228.3066275031833
$ws.Range("R3").Value = 2054.75964752865
$ws.Range("S3").Value = 0.02478330138626614
$ws.Range("T3").Value = 0.02478330138626614

# Row 4
$ws.Range("G4").Value = 32.93949833333333
$ws.Range("H4").Value = 98.818495
$ws.Range("I4").Value = 0.02571831923682078
$ws.Range("J4").Value = 0.02571831923682077
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03527466666666667
$ws.Range("N4").Value = 0.105824
$ws.Range("O4").Value = 0.004904310471994002
$ws.Range("P4").Value = 0.004904310471994003
$ws.Range("Q4").Value = 1.161929823875556
$ws.Range("R4").Value = 10.45736841488
$ws.Range("S4").Value = 0.0001261306223552249
$ws.Range("T4").Value = 0.0001261306223552249

# Row 5
$ws.Range("I5").Value = 0.4140443484779395
$ws.Range("J5").Value = 0.4140443484779395
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2262196666666666
$ws.Range("N5").Value = 0.6786589999999999
$ws.Range("O5").Value = 0.03145179203784564
$ws.Range("P5").Value = 0.03145179203784564
$ws.Range("Q5").Value = 119.9641875098899
$ws.Range("R5").Value = 1079.677687589009
$ws.Range("S5").Value = 0.01302243674277344
$ws.Range("T5").Value = 0.01302243674277344

# Row 6
$ws.Range("I6").Value = 0.4140443484779395
$ws.Range("J6").Value = 0.4140443484779395
$ws.Range("O6").Value = 0.9636438974901603
$ws.Range("P6").Value = 0.9636438974901604
$ws.Range("S6").Value = 0.3989913097010557
$ws.Range("T6").Value = 0.3989913097010558

# Row 7
$ws.Range("I7").Value = 0.4140443484779395
$ws.Range("J7").Value = 0.4140443484779395
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03527466666666667
$ws.Range("N7").Value = 0.105824
$ws.Range("O7").Value = 0.004904310471994002
$ws.Range("P7").Value = 0.004904310471994003
$ws.Range("Q7").Value = 18.70613987149155
$ws.Range("R7").Value = 168.355258843424
$ws.Range("S7").Value = 0.002030602034110293
$ws.Range("T7").Value = 0.002030602034110293

# Row 8
$ws.Range("G8").Value = 422.1807963333333
$ws.Range("H8").Value = 1266.542389
$ws.Range("I8").Value = 0.3296279860087694
$ws.Range("J8").Value = 0.3296279860087693
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2262196666666666
$ws.Range("N8").Value = 0.6786589999999999
$ws.Range("O8").Value = 0.03145179203784564
$ws.Range("P8").Value = 0.03145179203784564
$ws.Range("Q8").Value = 95.50559901959454
$ws.Range("R8").Value = 859.5503911763509
$ws.Range("S8").Value = 0.01036739086580171
$ws.Range("T8").Value = 0.01036739086580171

# Row 9
$ws.Range("G9").Value = 422.1807963333333
$ws.Range("H9").Value = 1266.542389
$ws.Range("I9").Value = 0.3296279860087694
$ws.Range("J9").Value = 0.3296279860087693
$ws.Range("O9").Value = 0.9636438974901603
$ws.Range("P9").Value = 0.9636438974901604
$ws.Range("Q9").Value = 2926.173095658003
$ws.Range("R9").Value = 26335.55786092203
$ws.Range("S9").Value = 0.3176439971593226
$ws.Range("T9").Value = 0.3176439971593225

# Row 10
$ws.Range("G10").Value = 422.1807963333333
$ws.Range("H10").Value = 1266.542389
$ws.Range("I10").Value = 0.3296279860087694
$ws.Range("J10").Value = 0.3296279860087693
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.03527466666666667
$ws.Range("N10").Value = 0.105824
$ws.Range("O10").Value = 0.004904310471994002
$ws.Range("P10").Value = 0.004904310471994003
$ws.Range("Q10").Value = 14.89228686372622
$ws.Range("R10").Value = 134.030581773536
$ws.Range("S10").Value = 0.0016165979836451
$ws.Range("T10").Value = 0.0016165979836451

# Row 11
$ws.Range("G11").Value = 16.509264
$ws.Range("H11").Value = 49.527792
$ws.Range("I11").Value = 0.01289001179132366
$ws.Range("J11").Value = 0.01289001179132366
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.2262196666666666
$ws.Range("N11").Value = 0.6786589999999999
$ws.Range("O11").Value = 0.03145179203784564
$ws.Range("P11").Value = 0.03145179203784564
$ws.Range("Q11").Value = 3.734720198991999
$ws.Range("R11").Value = 33.61248179092799
$ws.Range("S11").Value = 0.0004054139702260899
$ws.Range("T11").Value = 0.0004054139702260899

# Row 12
$ws.Range("G12").Value = 16.509264
$ws.Range("H12").Value = 49.527792
$ws.Range("I12").Value = 0.01289001179132366
$ws.Range("J12").Value = 0.01289001179132366
$ws.Range("O12").Value = 0.9636438974901603
$ws.Range("P12").Value = 0.9636438974901604
$ws.Range("Q12").Value = 114.42719461776
$ws.Range("R12").Value = 1029.84475155984
$ws.Range("S12").Value = 0.01242138120128525
$ws.Range("T12").Value = 0.01242138120128526

# Row 13
$ws.Range("G13").Value = 16.509264
$ws.Range("H13").Value = 49.527792
$ws.Range("I13").Value = 0.01289001179132366
$ws.Range("J13").Value = 0.01289001179132366
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.03527466666666667
$ws.Range("N13").Value = 0.105824
$ws.Range("O13").Value = 0.004904310471994002
$ws.Range("P13").Value = 0.004904310471994003
$ws.Range("Q13").Value = 0.582358784512
$ws.Range("R13").Value = 5.241229060608
$ws.Range("S13").Value = 0.00006321661981231479
$ws.Range("T13").Value = 0.00006321661981231481

# Row 14
$ws.Range("G14").Value = 44.62094166666666
$ws.Range("H14").Value = 133.862825
$ws.Range("I14").Value = 0.03483889192294087
$ws.Range("J14").Value = 0.03483889192294087
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.2262196666666666
$ws.Range("N14").Value = 0.6786589999999999
$ws.Range("O14").Value = 0.03145179203784564
$ws.Range("P14").Value = 0.03145179203784564
$ws.Range("Q14").Value = 10.09413455018611
$ws.Range("R14").Value = 90.84721095167498
$ws.Range("S14").Value = 0.001095745583589317
$ws.Range("T14").Value = 0.001095745583589317

# Row 15
$ws.Range("G15").Value = 44.62094166666666
$ws.Range("H15").Value = 133.862825
$ws.Range("I15").Value = 0.03483889192294087
$ws.Range("J15").Value = 0.03483889192294087
$ws.Range("O15").Value = 0.9636438974901603
$ws.Range("P15").Value = 0.9636438974901604
$ws.Range("Q15").Value = 309.2717625764166
$ws.Range("R15").Value = 2783.44586318775
$ws.Range("S15").Value = 0.03357228559686121
$ws.Range("T15").Value = 0.03357228559686121

# Row 16
$ws.Range("G16").Value = 44.62094166666666
$ws.Range("H16").Value = 133.862825
$ws.Range("I16").Value = 0.03483889192294087
$ws.Range("J16").Value = 0.03483889192294087
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03527466666666667
$ws.Range("N16").Value = 0.105824
$ws.Range("O16").Value = 0.004904310471994002
$ws.Range("P16").Value = 0.004904310471994003
$ws.Range("Q16").Value = 1.573988843644444
$ws.Range("R16").Value = 14.1658995928
$ws.Range("S16").Value = 0.0001708607424903462
$ws.Range("T16").Value = 0.0001708607424903462

# Row 17
$ws.Range("G17").Value = 234.229538
$ws.Range("H17").Value = 702.6886139999999
$ws.Range("I17").Value = 0.1828804425622059
$ws.Range("J17").Value = 0.1828804425622059
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.2262196666666666
$ws.Range("N17").Value = 0.6786589999999999
$ws.Range("O17").Value = 0.03145179203784564
$ws.Range("P17").Value = 0.03145179203784564
$ws.Range("Q17").Value = 52.98732800984732
$ws.Range("R17").Value = 476.8859520886259
$ws.Range("S17").Value = 0.005751917647255674
$ws.Range("T17").Value = 0.005751917647255674

# Row 18
$ws.Range("G18").Value = 234.229538
$ws.Range("H18").Value = 702.6886139999999
$ws.Range("I18").Value = 0.1828804425622059
$ws.Range("J18").Value = 0.1828804425622059
$ws.Range("O18").Value = 0.9636438974901603
$ws.Range("P18").Value = 0.9636438974901604
$ws.Range("Q18").Value = 1623.46600853642
$ws.Range("R18").Value = 14611.19407682778
$ws.Range("S18").Value = 0.1762316224453695
$ws.Range("T18").Value = 0.1762316224453695

# Row 19
$ws.Range("G19").Value = 234.229538
$ws.Range("H19").Value = 702.6886139999999
$ws.Range("I19").Value = 0.1828804425622059
$ws.Range("J19").Value = 0.1828804425622059
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0.3333333333333333
$ws.Range("M19").Value = 0.03527466666666667
$ws.Range("N19").Value = 0.105824
$ws.Range("O19").Value = 0.004904310471994002
$ws.Range("P19").Value = 0.004904310471994003
$ws.Range("Q19").Value = 8.262368876437332
$ws.Range("R19").Value = 74.36131988793599
$ws.Range("S19").Value = 0.0008969024695807239
$ws.Range("T19").Value = 0.000896902469580724
